$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ContractId values for stage rows 2-4 (BCD0017 -> BCD0013)
$ws.Range("B2").Value = "BCD0013"
$ws.Range("B3").Value = "BCD0013"
$ws.Range("B4").Value = "BCD0013"

# Update ContractId values for stage rows 5-7 (BCD0018 -> BCD0014)
$ws.Range("B5").Value = "BCD0014"
$ws.Range("B6").Value = "BCD0014"
$ws.Range("B7").Value = "BCD0014"

# Update TotalPaymentMade values (15000 -> 1000)
$ws.Range("N2").Value = 1000
$ws.Range("N3").Value = 1000
$ws.Range("N4").Value = 1000
$ws.Range("N5").Value = 1000
$ws.Range("N6").Value = 1000
$ws.Range("N7").Value = 1000

# Reflect the last active cell selection left in the sheet when it was saved
$ws.Range("G19").Select() | Out-Null

$wb.Save()
